# Refresh the TestInputs sheet's test data (previously hardcoded values)
# to reflect the latest getData() output:
#   - "Search String " (trailing space) -> "Search String"
#   - add a new "mango women" value in column C next to the search string

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestInputs")

$ws.Range("A1").Value = "Search String"
$ws.Range("C1").Value = "mango women"

# Match the saved selection state on the TestInputs sheet
$ws.Range("C12").Select() | Out-Null
